$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.179.29"
$ws.Range("E2").Value = "  +3.43%  "
$ws.Range("D3").Value = "2.456.92"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'322.49"
$ws.Range("E5").Value = "  +3.37%  "
$ws.Range("D6").Value = "'105.43"
$ws.Range("E6").Value = "  +3.26%  "
$ws.Range("E7").Value = "  +0.81%  "
$ws.Range("E9").Value = "  +5.28%  "
$ws.Range("D10").Value = "'36.15"
$ws.Range("E10").Value = "  +1.78%  "
$ws.Range("D11").Value = "'0.0808"
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("D13").Value = "'18.37"
$ws.Range("E13").Value = "  -2.35%  "
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("D15").Value = "2.843.16"
$ws.Range("E15").Value = "  +1.06%  "
$ws.Range("D16").Value = "2.446.02"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").Value = "'0.845"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").Value = "46.088.41"
$ws.Range("E18").Value = "  +3.42%  "
$ws.Range("D19").Value = "'12.62"
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("D21").Value = "0.0₃0940"
$ws.Range("E21").Value = "  +3.30%  "
$ws.Range("D22").Value = "'71.88"
$ws.Range("E22").Value = "  +4.24%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'248.21"
$ws.Range("E23").Value = "  +2.73%  "
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D24").Value = "'2.36"
$ws.Range("E24").Value = "  +1.71%  "
$ws.Range("E25").Value = "  +1.54%  "
$ws.Range("D26").Value = "'26.06"
$ws.Range("E26").Value = "  +2.92%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  -3.49%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("D30").Value = "'33.97"
$ws.Range("E30").Value = "  +1.31%  "
$ws.Range("D31").Value = "'49.31"
$ws.Range("E31").Value = "  +1.07%  "
$ws.Range("E32").Value = "  +6.02%  "
$ws.Range("D33").Value = "'20.39"
$ws.Range("E33").Value = "  +4.44%  "
$ws.Range("D34").Value = "'5.32"
$ws.Range("E34").Value = "  +2.01%  "
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("E37").Value = "  +0.90%  "
$ws.Range("E38").Value = "  -0.51%  "
$ws.Range("E39").Value = "  +0.36%  "
$ws.Range("D40").Value = "'128.43"
$ws.Range("E40").Value = "  +1.39%  "
$ws.Range("E41").Value = "  +3.45%  "
$ws.Range("E42").Value = "  +1.54%  "
$ws.Range("D43").Value = "'20.93"
$ws.Range("E43").Value = "  -4.72%  "
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("D45").Value = "1.960.15"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("E47").Value = "  -2.71%  "
$ws.Range("D48").Value = "'1.84"
$ws.Range("E48").Value = "  +9.31%  "
$ws.Range("D49").Value = "'9.21"
$ws.Range("E49").Value = "  -5.59%  "
$ws.Range("D50").Value = "'77.65"
$ws.Range("E50").Value = "  +4.93%  "
$ws.Range("D51").Value = "'4.91"
$ws.Range("E51").Value = "  +5.81%  "
